$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update USERID cell (G2): 33599 -> 44912
$ws.Range("G2").Value = 44912

# Update PREPERATION cell (F2): embedded user id 33599 -> 44912
$ws.Range("F2").Value = "Username : 44912,
Password : bni1234,
Cetak Laporan PDF,
Nama Laporan : Arus Kas,
Tipe Laporan : Konsolidasi Harian,
Produk : - ,
Mata Uang : IDR,
Status Posting : Posting ,
Tanggal Transaksi : 01/08/2022,
Tanggal Pembanding : 01/08/2022"

# Update sheet view: selected cell and scroll position (top-left visible cell = E1)
$ws.Range("J2").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
